# The car-sampling results were re-generated with an extra reporting-type
# row ("AW & ST: App") inserted into both the mean/IQR sampling-delay table
# and the cumulative sampling-delay table. Everything else in the workbook
# (the other sheets' "AW & ST: App" labels, etc.) automatically keeps
# resolving correctly because it's the same shared text value.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Sampling_Delay_Mean_IQR ---
# Insert a new row 6 for "AW & ST: App", pushing the existing
# "AW: Sensor, ..." rows down from 6/7/8 to 7/8/9.
$ws1 = $wb.Worksheets.Item("Sampling_Delay_Mean_IQR")
$ws1.Rows.Item(6).Insert()

$ws1.Cells.Item(6,1).Value = "AW & ST: App"
$ws1.Cells.Item(6,1).Font.Bold = $true
$ws1.Cells.Item(6,1).HorizontalAlignment = -4108
$ws1.Cells.Item(6,1).VerticalAlignment = -4160
$ws1.Cells.Item(6,1).Borders.LineStyle = 1

$ws1.Cells.Item(6,2).Value = 0.9399999999999999
$ws1.Cells.Item(6,3).Value = 1.18
$ws1.Cells.Item(6,4).Value = 1.49
$ws1.Cells.Item(6,5).Value = 1.51
$ws1.Cells.Item(6,6).Value = 2.03
$ws1.Cells.Item(6,7).Value = 1.78
$ws1.Cells.Item(6,8).Value = 2.47
$ws1.Cells.Item(6,9).Value = 2.05
$ws1.Cells.Item(6,10).Value = 2.92
$ws1.Cells.Item(6,11).Value = 2.15

# --- Sheet 2: Cumulative_Sampling_Delay ---
# Same new row for the cumulative-delay table.
$ws2 = $wb.Worksheets.Item("Cumulative_Sampling_Delay")
$ws2.Rows.Item(6).Insert()

$ws2.Cells.Item(6,1).Value = "AW & ST: App"
$ws2.Cells.Item(6,1).Font.Bold = $true
$ws2.Cells.Item(6,1).HorizontalAlignment = -4108
$ws2.Cells.Item(6,1).VerticalAlignment = -4160
$ws2.Cells.Item(6,1).Borders.LineStyle = 1

$ws2.Cells.Item(6,2).Value = 1.66
$ws2.Cells.Item(6,3).Value = 1.58

Write-Output "done"
